$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (row 1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values
$ws.Range("B2").Value = 61.478991029857156
$ws.Range("C2").Value = 73.205502253884291
$ws.Range("D2").Value = 57.138090646328997
$ws.Range("E2").Value = 78.809580296614044

# Update row 3 values
$ws.Range("B3").Value = 47.41038808743189
$ws.Range("C3").Value = 62.682346474154272
$ws.Range("D3").Value = 57.619929157500792
$ws.Range("E3").Value = 83.748390417755843

# Update the selection to match the new sqref B1:E3
$ws.Range("B1:E3").Select()
